$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NATMI LR-pairs table with newly-computed TPM-based values.
# Columns A-D (Sending cluster / Ligand / Receptor / Target cluster) and
# K-L (Receptor-expressing cells / Receptor detection rate) are unchanged.
# Only the TPM-dependent metric columns below are updated, keyed by
# 1-based column index: E=5 F=6 G=7 H=8 I=9 J=10 M=13 N=14 O=15 P=16
# Q=17 R=18 S=19 T=20
$data = @{
    2 = @{ 5=3.0; 6=1.0; 7=13.26223533333333; 8=39.786706; 9=0.391007655706778; 10=0.3910076557067781; 13=62.87391666666667; 14=188.62175; 15=0.5484251561826182; 16=0.5484251561826182; 17=833.8486791617221; 18=7504.638112455499; 19=0.2144384346495891; 20=0.2144384346495892 }
    3 = @{ 5=3.0; 6=1.0; 7=13.26223533333333; 8=39.786706; 9=0.391007655706778; 10=0.3910076557067781; 13=35.24551066666667; 14=105.736532; 15=0.3074331251635; 16=0.3074331251635; 17=467.4342569048435; 18=4206.908312143592; 19=0.1202087055567886; 20=0.1202087055567886 }
    4 = @{ 5=3.0; 6=1.0; 7=13.26223533333333; 8=39.786706; 9=0.391007655706778; 10=0.3910076557067781; 13=7.926563000000001; 14=23.779689; 15=0.06914038096772557; 16=0.06914038096772555; 17=105.1239438904927; 18=946.1154950144339; 19=0.0270344182768639; 20=0.0270344182768639 }
    5 = @{ 5=3.0; 6=1.0; 7=13.26223533333333; 8=39.786706; 9=0.391007655706778; 10=0.3910076557067781; 13=8.598489333333333; 14=25.795468; 15=0.07500133768615619; 16=0.07500133768615617; 17=114.0351890498231; 18=1026.316701448408; 19=0.02932609722353635; 20=0.02932609722353635 }
    6 = @{ 5=3.0; 6=1.0; 7=19.55844; 8=58.67532; 9=0.5766373150128344; 10=0.5766373150128344; 13=62.87391666666667; 14=188.62175; 15=0.5484251561826182; 16=0.5484251561826182; 17=1229.71572669; 18=11067.44154021; 19=0.3162424095466393; 20=0.3162424095466393 }
    7 = @{ 5=3.0; 6=1.0; 7=19.55844; 8=58.67532; 9=0.5766373150128344; 10=0.5766373150128344; 13=35.24551066666667; 14=105.736532; 15=0.3074331251635; 16=0.3074331251635; 17=689.3472056433601; 18=6204.124850790241; 19=0.1772774118402853; 20=0.1772774118402853 }
    8 = @{ 5=3.0; 6=1.0; 7=19.55844; 8=58.67532; 9=0.5766373150128344; 10=0.5766373150128344; 13=7.926563000000001; 14=23.779689; 15=0.06914038096772557; 16=0.06914038096772555; 17=155.03120684172; 18=1395.28086157548; 19=0.03986892364019375; 20=0.03986892364019375 }
    9 = @{ 5=3.0; 6=1.0; 7=19.55844; 8=58.67532; 9=0.5766373150128344; 10=0.5766373150128344; 13=8.598489333333333; 14=25.795468; 15=0.07500133768615619; 16=0.07500133768615617; 17=168.17303771664; 18=1513.55733944976; 19=0.04324856998571602; 20=0.04324856998571601 }
    10 = @{ 5=3; 6=1; 7=0.8919193333333334; 8=2.675758; 9=0.02629626747232247; 10=0.02629626747232247; 13=62.87391666666667; 14=188.62175; 15=0.5484251561826182; 16=0.5484251561826182; 17=56.07846183738889; 18=504.7061565365; 19=0.01442153459552835; 20=0.01442153459552835 }
    11 = @{ 5=3; 6=1; 7=0.8919193333333334; 8=2.675758; 9=0.02629626747232247; 10=0.02629626747232247; 13=35.24551066666667; 14=105.736532; 15=0.3074331251635; 16=0.3074331251635; 17=31.43615237680623; 18=282.925371391256; 19=0.008084343689151388; 20=0.008084343689151388 }
    12 = @{ 5=3; 6=1; 7=0.8919193333333334; 8=2.675758; 9=0.02629626747232247; 10=0.02629626747232247; 13=7.926563000000001; 14=23.779689; 15=0.06914038096772557; 16=0.06914038096772555; 17=7.069854786584668; 18=63.62869307926201; 19=0.001818133951065585; 20=0.001818133951065585 }
    13 = @{ 5=3; 6=1; 7=0.8919193333333334; 8=2.675758; 9=0.02629626747232247; 10=0.02629626747232247; 13=8.598489333333333; 14=25.795468; 15=0.07500133768615619; 16=0.07500133768615617; 17=7.669158873860445; 18=69.022429864744; 19=0.001972255236577142; 20=0.001972255236577142 }
    14 = @{ 5=2.0; 6=0.6666666666666666; 7=0.2055016666666667; 8=0.6165050000000001; 9=0.006058761808064916; 10=0.006058761808064917; 13=62.87391666666667; 14=188.62175; 15=0.5484251561826182; 16=0.5484251561826182; 17=12.92069466486111; 18=116.28625198375; 19=0.003322777390861284; 20=0.003322777390861285 }
    15 = @{ 5=2.0; 6=0.6666666666666666; 7=0.2055016666666667; 8=0.6165050000000001; 9=0.006058761808064916; 10=0.006058761808064917; 13=35.24551066666667; 14=105.736532; 15=0.3074331251635; 16=0.3074331251635; 17=7.243011184517779; 18=65.18710066066002; 19=0.001862664077274655; 20=0.001862664077274655 }
    16 = @{ 5=2.0; 6=0.6666666666666666; 7=0.2055016666666667; 8=0.6165050000000001; 9=0.006058761808064916; 10=0.006058761808064917; 13=7.926563000000001; 14=23.779689; 15=0.06914038096772557; 16=0.06914038096772555; 17=1.628921907438334; 18=14.660297166945; 19=0.0004189050996023141; 20=0.000418905099602314 }
    17 = @{ 5=2.0; 6=0.6666666666666666; 7=0.2055016666666667; 8=0.6165050000000001; 9=0.006058761808064916; 10=0.006058761808064917; 13=8.598489333333333; 14=25.795468; 15=0.07500133768615619; 16=0.07500133768615617; 17=1.767003888815556; 18=15.90303499934; 19=0.000454415240326663; 20=0.000454415240326663 }
}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $ws.Cells.Item($r, $c).Value = $data[$r][$c]
    }
}